$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US".
#    This string is shared by the Overview sheet (zh-cn/de-de status
#    columns E/F) and by the per-locale "Status" column (C) on the zh-cn and
#    de-de sheets - update every cell that shows it so the shared string
#    itself is rewritten (not just re-pointed).
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("C3").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# 2. zh-cn sheet: the handback run has produced a target file and a handback
#    file for each row - fill in "Latest Target File" (I) / "Latest Handback
#    File" (J), with the target-file cell hyperlinked like column A.
# ---------------------------------------------------------------------------
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/345667f0f8bc8cccb2f29798252e79e750655bde/e2e/27a7ff87-1504-46fb-a672-5bd3229e647d.md", "", "", "27a7ff87-1504-46fb-a672-5bd3229e647d.md")
$zhcn.Range("J2").Value = "27a7ff87-1504-46fb-a672-5bd3229e647d.208fe2b2f3d3998e3453a283b54a07095b20552b.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-10-18 12:27:03"

$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/345667f0f8bc8cccb2f29798252e79e750655bde/e2e/79b83687-3609-4f18-b3c8-7a982e7f2368.md", "", "", "79b83687-3609-4f18-b3c8-7a982e7f2368.md")
$zhcn.Range("J3").Value = "79b83687-3609-4f18-b3c8-7a982e7f2368.a6eba2a658119e831227b2be51bdb917fec73164.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-10-18 12:27:03"

$zhcn.Columns.Item(3).ColumnWidth = 29.9777050018311
$zhcn.Columns.Item(9).ColumnWidth = 40
$zhcn.Columns.Item(10).ColumnWidth = 40

# ---------------------------------------------------------------------------
# 3. de-de sheet: same as zh-cn, plus the handback completion timestamp (K).
# ---------------------------------------------------------------------------
$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/345667f0f8bc8cccb2f29798252e79e750655bde/e2e/27a7ff87-1504-46fb-a672-5bd3229e647d.md", "", "", "27a7ff87-1504-46fb-a672-5bd3229e647d.md")
$dede.Range("J2").Value = "27a7ff87-1504-46fb-a672-5bd3229e647d.208fe2b2f3d3998e3453a283b54a07095b20552b.de-de.xlf"
$dede.Range("K2").Value = "2016-10-18 12:27:20"

$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/345667f0f8bc8cccb2f29798252e79e750655bde/e2e/79b83687-3609-4f18-b3c8-7a982e7f2368.md", "", "", "79b83687-3609-4f18-b3c8-7a982e7f2368.md")
$dede.Range("J3").Value = "79b83687-3609-4f18-b3c8-7a982e7f2368.a6eba2a658119e831227b2be51bdb917fec73164.de-de.xlf"
$dede.Range("K3").Value = "2016-10-18 12:27:20"

$dede.Columns.Item(3).ColumnWidth = 29.9777050018311
$dede.Columns.Item(9).ColumnWidth = 40
$dede.Columns.Item(10).ColumnWidth = 40
